# "Generate Report for Handoff"
#
# A new handoff pass was run for the cfb95521-ef8c-41f7-9885-bb87de814fc9
# file (the last row, row 7, on every sheet). This refreshes the
# "Latest Handoff" timestamps that the report shows for that row:
#   - Overview!G7            "Latest HO Xliff Generate Date"
#   - zh-cn!H7                "Latest Handoff Datetime"
#   - de-de!H7                "Latest Handoff Datetime"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-13 04:50:23"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-13 04:50:16"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-13 04:50:23"
